$d = $word.ActiveDocument
$d.Content.Find.Execute("Arbitrary Left Back", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Arbitrary Back", 2)
